$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.919
$ws.Range("B2").Value = 0.979
$ws.Range("C2").Value = 0.934
$ws.Range("D2").Value = 0.923
$ws.Range("E2").Value = 0.944
$ws.Range("F2").Value = 0.92
